$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 28000000.0
$ws.Range("C4").Value = 27000000.0
$ws.Range("D4").Value = 26000000.0
$ws.Range("E4").Value = 26000000.0
$ws.Range("F4").Value = 23000000.0

# Row 12 - Accounts Payable
$ws.Range("B12").Value = 53000000.0
$ws.Range("C12").Value = 54000000.0
$ws.Range("D12").Value = 42000000.0
$ws.Range("E12").Value = 39000000.0
$ws.Range("F12").Value = 40000000.0

# Row 37 - Net Debt
$ws.Range("G37").Value = -173761000.0

# Row 38 - Total Debt
$ws.Range("G38").Value = 154424000.0
